# Update Data From Excel
# Re-assigns the "Ma_Mon_Hoc" (column D) course-code values for every student
# row on the active sheet. The codes now repeat in contiguous blocks of 7 rows
# (one course code per class group) instead of the old round-robin pattern,
# and the handful of course codes that are no longer used by anyone
# (0804-01, 0806-01, 0807-01) are dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new Ma_Mon_Hoc value, grouped in blocks of 7 rows each.
$maMonHoc = @{
    2 = "0801-01";  3 = "0801-01";  4 = "0801-01";  5 = "0801-01"
    6 = "0801-01";  7 = "0801-01";  8 = "0801-01"

    9 = "0805-01"; 10 = "0805-01"; 11 = "0805-01"; 12 = "0805-01"
    13 = "0805-01"; 14 = "0805-01"; 15 = "0805-01"

    16 = "0803-01"; 17 = "0803-01"; 18 = "0803-01"; 19 = "0803-01"
    20 = "0803-01"; 21 = "0803-01"; 22 = "0803-01"

    23 = "0808-01"; 24 = "0808-01"; 25 = "0808-01"; 26 = "0808-01"
    27 = "0808-01"; 28 = "0808-01"; 29 = "0808-01"

    30 = "0802-01"; 31 = "0802-01"; 32 = "0802-01"; 33 = "0802-01"
    34 = "0802-01"; 35 = "0802-01"; 36 = "0802-01"
}

foreach ($row in 2..36) {
    $ws.Cells.Item($row, 4).Value = $maMonHoc[$row]
}

# Match the saved selection / active cell from the refreshed data pull.
$ws.Range("D2:D36").Select()

# The workbook was re-saved with a portrait page setup.
$ws.PageSetup.Orientation = 1
